$wb = $excel.ActiveWorkbook

# Overview!G2 and de-de!H2 share the same original text
# ("2016-08-26 01:04:10") and both move to the new timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 01:05:13"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 01:05:05"
$wsZhCn.Range("K2").Value = "2016-08-26 01:05:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-26 01:05:13"
$wsDeDe.Range("K2").Value = "2016-08-26 01:05:36"
